# Manual_Speciation_Calculation.xlsx - split fugacity calc into chunks,
# add a "Speciation" worksheet, and re-sort the "gammas" worksheet by species.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# 1. Re-sort "gammas" sheet rows into (O2, H2, S2, CO, CO2, SO2, H2S, H2O)
#    order -- same label/value pairs, just re-ordered rows 2-9.
# ---------------------------------------------------------------------------
$gammas = $wb.Worksheets.Item("gammas")

$gammaOrder = @(
    @("O2",  1.19946588566002),
    @("H2",  1.1995158606479099),
    @("S2",  1.16052802963702),
    @("CO",  1.2702307224970899),
    @("CO2", 1.17702986434571),
    @("SO2", 1.12402533093521),
    @("H2S", 1.11989805891392),
    @("H2O", 0.86908487374434096)
)

for ($i = 0; $i -lt $gammaOrder.Count; $i++) {
    $row = $i + 2
    $gammas.Cells.Item($row, 1).Value = $gammaOrder[$i][0]
    $gammas.Cells.Item($row, 2).Value = $gammaOrder[$i][1]
}

$gammas.Range("C26").Select()

# ---------------------------------------------------------------------------
# 2. "fugacities" sheet: split the E column "difference" formulas into a
#    shared formula block, restyle D/E columns, and clear the explicit
#    column styles that are no longer used.
# ---------------------------------------------------------------------------
$fug = $wb.Worksheets.Item("fugacities")

$fug.Range("E2:E9").Formula = "=C2-D2"

$fug.Columns.Item(4).ClearFormats()
$fug.Columns.Item(5).ClearFormats()

$fug.Range("D2").NumberFormat = "0.00E+00"
$fug.Range("D2").Font.Name = "Andale Mono"
$fug.Range("D2").Font.Size = 12

$fug.Range("D3:D9").Font.Name = "Andale Mono"
$fug.Range("D3:D9").Font.Size = 12

$fug.Range("E2:E9").NumberFormat = "0.00E+00"
$fug.Range("E2:E9").Font.Name = "Andale Mono"
$fug.Range("E2:E9").Font.Size = 12

$fug.TabSelected = $false
$fug.Range("B28").Select()

# ---------------------------------------------------------------------------
# 3. Add the new "Speciation" sheet (after "fugacities", i.e. last tab).
# ---------------------------------------------------------------------------
$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$spec = $wb.Worksheets.Add([System.Reflection.Missing]::Value, $lastSheet)
$spec.Name = "Speciation"

$spec.Columns("B:E").ColumnWidth = 12.1640625

$spec.Range("B1").Value = "mole fraction"
$spec.Range("C1").Value = "norm_fix_ox"
$spec.Range("D1").Value = "MW"
$spec.Range("E1").Value = "MPO"
$spec.Range("F1").Value = "wt%"
$spec.Range("B1:F1").Font.Bold = $true

$specOrder = @("O2", "H2", "S2", "CO", "CO2", "SO2", "H2S", "H2O")
for ($i = 0; $i -lt $specOrder.Count; $i++) {
    $row = $i + 2
    $spec.Cells.Item($row, 1).Value = $specOrder[$i]
}

$spec.Range("B2").Formula = "=fugacities!C2/(gammas!B2*Sheet1!`$B`$9)"
$spec.Range("B3:B9").Formula = "=fugacities!C3/(gammas!B3*Sheet1!`$B`$9)"
$spec.Range("B3").Formula = "=fugacities!C3/(gammas!B3*Sheet1!`$B`$9)"
$spec.Range("B4").Formula = "=fugacities!C4/(gammas!B4*Sheet1!`$B`$9)"
$spec.Range("B5").Formula = "=fugacities!C5/(gammas!B5*Sheet1!`$B`$9)"
$spec.Range("B6").Formula = "=fugacities!C6/(gammas!B6*Sheet1!`$B`$9)"
$spec.Range("B7").Formula = "=fugacities!C7/(gammas!B7*Sheet1!`$B`$9)"
$spec.Range("B8").Formula = "=fugacities!C8/(gammas!B8*Sheet1!`$B`$9)"
$spec.Range("B9").Formula = "=fugacities!C9/(gammas!B9*Sheet1!`$B`$9)"

$spec.Range("C2").Formula = "=B2"
$spec.Range("C3:C9").Formula = "=B3/B`$12*(1-B`$2)"

$spec.Range("D2").Formula = "=2*15.999"
$spec.Range("D3").Formula = "=2.016"
$spec.Range("D4").Formula = "=2*32.065"
$spec.Range("D5").Formula = "=20.01"
$spec.Range("D6").Value = 44.01
$spec.Range("D7").Value = 64.066
$spec.Range("D8").Value = 34.1
$spec.Range("D9").Value = 18.015

$spec.Range("E2").Formula = "=C2*D2"
$spec.Range("E3:E9").Formula = "=C3*D3"

$spec.Range("F2").Formula = "=100*E2/E`$11"
$spec.Range("F3:F9").Formula = "=100*E3/E`$11"

$spec.Range("A11").Value = "sum"
$spec.Range("B11").Formula = "=SUM(B2:B9)"
$spec.Range("C11").Formula = "=SUM(C2:C9)"
$spec.Range("E11").Formula = "=SUM(E2:E9)"
$spec.Range("F11").Formula = "=SUM(F2:F9)"

$spec.Range("A12").Value = "normsum"
$spec.Range("B12").Formula = "=SUM(B3:B9)"

$spec.Range("F2").Select()
$spec.Activate()
